$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10 of data (climate emergency gets added to shared strings first)
$ws.Range("A10").Value = "climate emergency"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 5

# Change header D1 from "normality" to "sign" (sign gets added next)
$ws.Range("D1").Value = "sign"

# word of the year gets added last
$ws.Range("E10").Value = "word of the year"

# Update selection to match final state
$ws.Range("D11").Select()
